$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.497.44"
$ws.Range("E2").Value = "  +2.31%  "

$ws.Range("D3").Value = "1.872.03"
$ws.Range("E3").Value = "  +1.63%  "

$ws.Range("D4").Formula = "'1.015"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.92%  "

$ws.Range("D5").Formula = "'312.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.21%  "

$ws.Range("D6").Formula = "'1.014"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.87%  "

$ws.Range("D7").Formula = "'0.4785"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.76%  "

$ws.Range("D8").Formula = "'0.3784"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.14%  "

$ws.Range("D9").Formula = "'0.07376"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.52%  "

$ws.Range("D10").Formula = "'0.9394"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.51%  "

$ws.Range("D11").Formula = "'20.74"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.72%  "

$ws.Range("E12").Value = "  +2.29%  "

$ws.Range("D13").Value = "1.892.60"
$ws.Range("E13").Value = "  -0.67%  "

$ws.Range("D14").Formula = "'5.449"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.67%  "

$ws.Range("D15").Formula = "'6.585"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.07%  "

$ws.Range("D16").Formula = "'91.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.81%  "

$ws.Range("D17").Formula = "'1.016"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.77%  "

$ws.Range("D18").Formula = "'0.000008924"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.47%  "


$ws.Range("D20").Formula = "'14.92"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.63%  "

$ws.Range("D21").Value = "27.538.64"
$ws.Range("E21").Value = "  +2.29%  "

$ws.Range("D22").Formula = "'5.137"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.77%  "

$ws.Range("D23").Formula = "'10.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.08%  "

$ws.Range("D24").Formula = "'1.964"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.32%  "

$ws.Range("D25").Formula = "'154.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.40%  "

$ws.Range("D26").Formula = "'18.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.41%  "

$ws.Range("D27").Formula = "'2.021"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.03%  "

$ws.Range("D28").Formula = "'116.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.58%  "

$ws.Range("D29").Formula = "'5.001"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.29%  "

$ws.Range("D30").Formula = "'0.08937"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.06%  "

$ws.Range("D31").Formula = "'3.339"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.67%  "

$ws.Range("D32").Formula = "'1.221"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.33%  "

$ws.Range("D33").Formula = "'4.615"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.09%  "

$ws.Range("D34").Formula = "'0.7539"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("D35").Formula = "'2.713"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.04%  "

$ws.Range("D36").Formula = "'0.02063"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.87%  "

$ws.Range("D37").Formula = "'1.118"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.76%  "

$ws.Range("D38").Formula = "'0.05306"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.99%  "

$ws.Range("D39").Formula = "'3.003"
$ws.Range("D39").Style = "Normal"

$ws.Range("D40").Formula = "'0.5366"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.11%  "

$ws.Range("D41").Formula = "'7.085"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.68%  "

$ws.Range("D42").Formula = "'0.1530"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.22%  "

$ws.Range("D43").Formula = "'8.437"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.72%  "

$ws.Range("D44").Formula = "'10.60"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.16%  "

$ws.Range("D45").Formula = "'0.4826"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.15%  "

$ws.Range("E46").Value = "  +0.92%  "

$ws.Range("D47").Formula = "'1.664"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.94%  "

$ws.Range("D48").Formula = "'102.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.16%  "

$ws.Range("D49").Formula = "'67.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.68%  "

$ws.Range("D50").Formula = "'0.06099"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.27%  "

$ws.Range("D51").Formula = "'0.9342"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.55%  "
